# Automatic update of files.
#
# 1) Bump the "Förändrad" (changed/updated) date in column C for every
#    data row (2..61) from 45184 (2023-09-15) to 45186 (2023-09-17).
# 2) Add a friendly-text second argument to every HYPERLINK() formula
#    (columns S, T, V, W, X, Y) so the link label shows the case id
#    from column A instead of the raw URL.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 61
$newDate = 45186

# --- 1) Update the "Förändrad" date column (C) for all data rows ---
for ($r = 2; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    if ($cCell.Formula() -ne "") {
        $cCell.Value = $newDate
    }
}

# --- 2) Add the case-id label as the second HYPERLINK() argument ---
$linkCols = @("S", "T", "V", "W", "X", "Y")

for ($r = 2; $r -le $lastRow; $r++) {
    $label = $ws.Range("A$r").Value()
    if ([string]::IsNullOrEmpty($label)) {
        continue
    }

    foreach ($col in $linkCols) {
        $cell = $ws.Range("$col$r")
        $formula = $cell.Formula()

        if ($formula -and $formula.ToUpper().Contains("HYPERLINK(") -and -not $formula.Contains(",")) {
            # Formula looks like: =HYPERLINK("url")  -> add `, "label"` before closing paren
            $trimmed = $formula.TrimEnd()
            $newFormula = $trimmed.Substring(0, $trimmed.Length - 1) + ', "' + $label + '")'
            $cell.Formula = $newFormula
        }
    }
}
